$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 (Subj) values for B1:E1
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 (CON) values for B2:E2
$ws.Range("B2").Value = 16.823171312272251
$ws.Range("C2").Value = 14.511828908386349
$ws.Range("D2").Value = 13.483211846899906
$ws.Range("E2").Value = 1.1081228576872775

# Update row 3 (STR) values for B3:E3
$ws.Range("B3").Value = 34.135159043676282
$ws.Range("C3").Value = 5.2032961379966878
$ws.Range("D3").Value = 3.6095552557492852
$ws.Range("E3").Value = 3.281647899050256

# Update the selection to match the new selected range
$ws.Range("B1:E3").Select()
